# Applies the commit's changes to CATALOGO_GRUPOS.xlsx:
#  - Adds two new product codes (evol1000, PERFA0261) near the top of both
#    the "MEJORAR" and "PREMIUM" sheets' lists, pushing the existing rows
#    down.
#  - Switches the active sheet/tab back to "MEJORAR" (was "PREMIUM").
#  - Updates the selections left behind on both sheets.

$wb = $excel.ActiveWorkbook

$xlCenter = -4108
$xlLeft   = -4131

# ---------------------------------------------------------------------
# Sheet "MEJORAR"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("MEJORAR")

$mejorarValues = @(
    "Codigo_Producto",
    "evol1000",
    "evol0330",
    "PERFA0261",
    "evol0088",
    "evol0330",
    "evol3089",
    "evol0025",
    "evol3245",
    "evol1970"
)

for ($i = 0; $i -lt $mejorarValues.Length; $i++) {
    $row = $i + 1
    $cell = $ws1.Cells.Item($row, 1)
    $cell.Value = $mejorarValues[$i]
    if ($row -ne 1) {
        $cell.HorizontalAlignment = $xlCenter
    }
}

# ---------------------------------------------------------------------
# Sheet "PREMIUM"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PREMIUM")

# Row 2 used to hold "evol0088" styled with left-alignment (s=3); that
# formatting now belongs to row 5 (where "evol0088" ends up), so the old
# cell's leftover formatting must be cleared before the new value goes in.
$ws2.Cells.Item(2, 1).ClearFormats()

$premiumValues = @(
    "Codigo_Producto",
    "evol1000",
    "evol0330",
    "PERFA0261",
    "evol0088",
    "evo115tu",
    "evol0028",
    "evol3510",
    "evorieg153",
    "evol0070",
    "evol2530",
    "evol0107",
    "evol0435"
)

for ($i = 0; $i -lt $premiumValues.Length; $i++) {
    $row = $i + 1
    $cell = $ws2.Cells.Item($row, 1)
    $cell.Value = $premiumValues[$i]
    if ($row -eq 5) {
        $cell.HorizontalAlignment = $xlLeft
    }
}

$ws2.Range("A2:A4").Select() | Out-Null

# ---------------------------------------------------------------------
# Re-activate "MEJORAR" as the selected tab (it was "PREMIUM" before) and
# leave its A2:A4 selection current, matching the saved view state.
# ---------------------------------------------------------------------
$ws1.Range("A2:A4").Select() | Out-Null
